$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure the Price column keeps its exact text representation
# (e.g. "1.00", "6.10", "523.61") instead of being auto-converted to a number.
$ws.Range("D2:D51").NumberFormat = "@"

# Data rows 2-51: Coin, Link, Price, Volume(1h)
$data = @(
    ,('Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '57.466.84', '  -0.73%  ')
    ,('Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '3.107.32', '  +1.17%  ')
    ,('TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  +0.02%  ')
    ,('BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '523.61', '  +1.33%  ')
    ,('Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '141.29', '  -0.84%  ')
    ,('USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.00', '  -0.01%  ')
    ,('LidoStakedEther', 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth', '3.106.91', '  +1.20%  ')
    ,('XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.437', '  +0.11%  ')
    ,('Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '7.21', '  -1.21%  ')
    ,('Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.108', '  +0.51%  ')
    ,('Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.386', '  +1.82%  ')
    ,('WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '3.644.53', '  +1.31%  ')
    ,('TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.131', '  +1.04%  ')
    ,('Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '26.11', '  -0.36%  ')
    ,('ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0000164', '  -0.14%  ')
    ,('WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '57.553.05', '  -0.57%  ')
    ,('WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '3.102.93', '  +1.19%  ')
    ,('Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '6.10', '  +0.37%  ')
    ,('Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '12.79', '  -0.79%  ')
    ,('Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '8.07', '  -0.87%  ')
    ,('BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '338.17', '  +1.85%  ')
    ,('Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.999', '  -0.30%  ')
    ,('Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.513', '  +2.63%  ')
    ,('Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '66.77', '  +1.72%  ')
    ,('Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.168', '  -0.66%  ')
    ,('Binance-PegBSC-USD', 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd', '1.00', '  +0.08%  ')
    ,('PEPE', 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', '0.0₃0913', '  +0.61%  ')
    ,('RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '6.51', '  +1.23%  ')
    ,('USDe', 'https://coinranking.com/coin/exbfr2U-0+usde-usde', '0.999', '  -0.02%  ')
    ,('InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '7.19', '  -0.59%  ')
    ,('PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.87', '  +2.53%  ')
    ,('EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '20.98', '  +1.12%  ')
    ,('Fetch.AI', 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet', '1.19', '  -0.27%  ')
    ,('Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '158.32', '  +2.33%  ')
    ,('NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '4.63', '  +1.91%  ')
    ,('Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '6.11', '  +1.92%  ')
    ,('EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '27.09', '  -0.15%  ')
    ,('ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '1.28', '  +1.19%  ')
    ,('Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.0659', '  -2.29%  ')
    ,('Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '3.95', '  +0.83%  ')
    ,('RenzoRestakedETH', 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth', '3.149.19', '  +1.15%  ')
    ,('Mantle', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt', '0.687', '  +4.55%  ')
    ,('Stacks', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', '1.52', '  +10.65%  ')
    ,('OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '36.84', '  +0.72%  ')
    ,('FirstDigitalUSD', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', '1.00', '  -0.04%  ')
    ,('Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '2.305.63', '  +1.79%  ')
    ,('VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.0259', '  +0.42%  ')
    ,('ONDO', 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo', '0.975', '  +3.67%  ')
    ,('InjectiveProtocol', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj', '20.78', '  -0.32%  ')
    ,('Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '6.01', '  +1.63%  ')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
}
